# Update the TPM-derived metrics for the Flt3l-Flt3 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.46387
$ws.Range("H2").Value = 25.39161
$ws.Range("I2").Value = 0.3683986519988784
$ws.Range("J2").Value = 0.3683986519988785
$ws.Range("M2").Value = 1.419902
$ws.Range("N2").Value = 4.259706
$ws.Range("Q2").Value = 12.01786594074
$ws.Range("R2").Value = 108.16079346666
$ws.Range("S2").Value = 0.3683986519988784
$ws.Range("T2").Value = 0.3683986519988785

# Row 3
$ws.Range("I3").Value = 0.3815382171230672
$ws.Range("J3").Value = 0.3815382171230673
$ws.Range("M3").Value = 1.419902
$ws.Range("N3").Value = 4.259706
$ws.Range("Q3").Value = 12.44650358999666
$ws.Range("R3").Value = 112.01853230997
$ws.Range("S3").Value = 0.3815382171230672
$ws.Range("T3").Value = 0.3815382171230673

# Row 4
$ws.Range("G4").Value = 3.186597333333334
$ws.Range("H4").Value = 9.559792000000002
$ws.Range("I4").Value = 0.138699928290867
$ws.Range("J4").Value = 0.1386999282908671
$ws.Range("M4").Value = 1.419902
$ws.Range("N4").Value = 4.259706
$ws.Range("Q4").Value = 4.524655926794667
$ws.Range("R4").Value = 40.721903341152
$ws.Range("S4").Value = 0.138699928290867
$ws.Range("T4").Value = 0.1386999282908671

# Row 5
$ws.Range("G5").Value = 2.558542666666666
$ws.Range("H5").Value = 7.675628
$ws.Range("I5").Value = 0.1113632025871871
$ws.Range("J5").Value = 0.1113632025871872
$ws.Range("M5").Value = 1.419902
$ws.Range("N5").Value = 4.259706
$ws.Range("Q5").Value = 3.632879849485332
$ws.Range("R5").Value = 32.695918645368
$ws.Range("S5").Value = 0.1113632025871871
$ws.Range("T5").Value = 0.1113632025871872
